# Apply updates described by the diff:
#  - tweak a handful of odds values in rows 2-4
#  - insert a brand new match row at row 5 (pushing the former row 5 down to row 6)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: small odds corrections -----------------------------------------
$ws.Range("Q2").Value = 1.67
$ws.Range("R2").Value = 2.15

# --- Row 3: small odds corrections -----------------------------------------
$ws.Range("G3").Value  = 2.1
$ws.Range("H3").Value  = 3.25
$ws.Range("I3").Value  = 3.1
$ws.Range("M3").Value  = 1.03
$ws.Range("N3").Value  = 9.5
$ws.Range("AI3").Value = 15
$ws.Range("AL3").Value = 26
$ws.Range("AN3").Value = 4.33
$ws.Range("AR3").Value = 67
$ws.Range("AW3").Value = 5

# --- Row 4: small odds corrections -----------------------------------------
$ws.Range("Q4").Value = 1.95
$ws.Range("R4").Value = 1.85

# --- Insert a brand new row at position 5 -----------------------------------
# This shifts the existing row 5 (xjgVeHTp / Asan - Daegu) down to row 6,
# unchanged, and leaves a blank row 5 ready to be filled with the new match.
$ws.Rows.Item(5).Insert()

$row5 = New-Object 'object[,]' 1,56
$row5[0,0]  = "Gzc7QLHb"
$row5[0,1]  = "28/11/2024"
$row5[0,2]  = "14:00"
$row5[0,3]  = "SAUDI ARABIA - SAUDI PROFESSIONAL LEAGUE"
$row5[0,4]  = "Al Fateh"
$row5[0,5]  = "Al Riyadh"
$row5[0,6]  = 2.35
$row5[0,7]  = 3.2
$row5[0,8]  = 2.8
$row5[0,9]  = 3
$row5[0,10] = 2.1
$row5[0,11] = 3.4
$row5[0,12] = 1.05
$row5[0,13] = 8.5
$row5[0,14] = 1.29
$row5[0,15] = 3.5
$row5[0,16] = 1.93
$row5[0,17] = 1.88
$row5[0,18] = 1.4
$row5[0,19] = 2.75
$row5[0,20] = 1.73
$row5[0,21] = 2
$row5[0,22] = 8.5
$row5[0,23] = 12
$row5[0,24] = 10
$row5[0,25] = 23
$row5[0,26] = 19
$row5[0,27] = 29
$row5[0,28] = 10
$row5[0,29] = 6.5
$row5[0,30] = 13
$row5[0,31] = 41
$row5[0,32] = 151
$row5[0,33] = 10
$row5[0,34] = 15
$row5[0,35] = 11
$row5[0,36] = 29
$row5[0,37] = 23
$row5[0,38] = 34
$row5[0,39] = 4.5
$row5[0,40] = 13
$row5[0,41] = 23
$row5[0,42] = 41
$row5[0,43] = 67
$row5[0,44] = 151
$row5[0,45] = 2.75
$row5[0,46] = 8
$row5[0,47] = 51
$row5[0,48] = 5
$row5[0,49] = 15
$row5[0,50] = 26
$row5[0,51] = 51
$row5[0,52] = 67
$row5[0,53] = 151
$row5[0,54] = 81
$row5[0,55] = 81

$ws.Range("A5:BD5").Value = $row5

Write-Output "done"
